$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header labels for the two new trailing columns (M, N) ---
$ws.Range("M1").Value = "V_hi (V p.u.)"
$ws.Range("N1").Value = "V_lo (V p.u.)"

# --- Updated timestamps (col B) -- this snapshot is from a later hunting run ---
$ws.Range("B2").Value = 44560.713970717603
$ws.Range("B3").Value = 44560.716146851853
$ws.Range("B4").Value = 44560.716822361108
$ws.Range("B5").Value = 44560.718002013891
$ws.Range("B6").Value = 44560.71876945602
$ws.Range("B7").Value = 44560.719873564813
$ws.Range("B8").Value = 44560.720654606477
$ws.Range("B9").Value = 44560.721680410657

# --- P_per_hi / Q_per_hi (cols H, I) changed for rows 2-3 in the new run ---
$ws.Range("H2").Value = -1279.2671069477101
$ws.Range("I2").Value = -619.16527976269151
$ws.Range("H3").Value = -1279.267106947708
$ws.Range("I3").Value = -619.16527976269072

# --- P_per_lo / Q_per_lo (cols K, L) changed for row 3 in the new run ---
$ws.Range("K3").Value = 736.88191505216741
$ws.Range("L3").Value = 356.650846885249

# --- New V_hi (V p.u.) / V_lo (V p.u.) values for every data row ---
$vhi = @{2=1.1618999999999999; 3=1.0765; 4=1.0936999999999999; 5=1.0431999999999999; 6=1.0833999999999999; 7=1.044; 8=1.1395; 9=1.071}
$vlo = @{2=1.0115000000000001; 3=0.94199999999999995; 4=0.99909999999999999; 5=0.94410000000000005; 6=0.9879; 7=0.94520000000000004; 8=0.99850000000000005; 9=0.93779999999999997}

foreach ($r in 2..9) {
    $ws.Cells.Item($r, 13).Value = $vhi[$r]
    $ws.Cells.Item($r, 14).Value = $vlo[$r]
}

# Copy the header style (bold, centered, bordered) from the last existing header (L1)
# onto the two new header cells so they look consistent with the rest of the row.
$ws.Range("L1").Copy() | Out-Null
$ws.Range("M1:N1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Match the column width used by the neighbouring numeric columns.
$ws.Columns("M:N").ColumnWidth = 12.83

# Move the selection to where the user left off after entering the new data.
$ws.Range("D10").Select() | Out-Null
